$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 111, "face/face072.png", "passen", "face"),
    @(3, 68, "flower/flower099.png", "heben", "flower"),
    @(4, 5, "flower/flower089.png", "bauen", "flower"),
    @(5, 40, "face/face071.png", "parken", "face"),
    @(6, 38, "face/face094.png", "kennen", "face"),
    @(7, 104, "face/face099.png", "betteln", "face"),
    @(8, 11, "face/face066.png", "kranken", "face"),
    @(9, 119, "flower/flower122.png", "küssen", "flower"),
    @(10, 19, "flower/flower075.png", "planen", "flower"),
    @(11, 45, "flower/flower071.png", "deuten", "flower"),
    @(12, 90, "flower/flower107.png", "zielen", "flower"),
    @(13, 17, "flower/flower078.png", "lügen", "flower"),
    @(14, 84, "face/face101.png", "fließen", "face"),
    @(15, 13, "flower/flower066.png", "holen", "flower"),
    @(16, 51, "face/face064.png", "sparen", "face"),
    @(17, 46, "flower/flower105.png", "quellen", "flower"),
    @(18, 30, "face/face097.png", "tollen", "face"),
    @(19, 28, "flower/flower081.png", "ändern", "flower"),
    @(20, 106, "face/face088.png", "nullen", "face"),
    @(21, 93, "flower/flower101.png", "trotzen", "flower"),
    @(22, 57, "flower/flower092.png", "reisen", "flower"),
    @(23, 118, "face/face119.png", "grenzen", "face"),
    @(24, 41, "face/face068.png", "prüfen", "face"),
    @(25, 73, "face/face102.png", "piepen", "face"),
    @(26, 70, "face/face096.png", "hassen", "face"),
    @(27, 122, "flower/flower067.png", "wählen", "flower"),
    @(28, 66, "flower/flower072.png", "heißen", "flower"),
    @(29, 113, "face/face103.png", "achten", "face"),
    @(30, 105, "flower/flower110.png", "öffnen", "flower"),
    @(31, 98, "face/face105.png", "ärgern", "face"),
    @(32, 0, "face/face092.png", "stoßen", "face"),
    @(33, 10, "flower/flower065.png", "münzen", "flower")
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Range("B$r").Value = $row[1]
    $ws.Range("C$r").Value = $row[2]
    $ws.Range("D$r").Value = $row[3]
    $ws.Range("E$r").Value = $row[4]
}